$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before the current row 10 (old row10 -> new row11, etc.)
$ws.Rows.Item(10).Insert()

$ws.Range("A10").Value = 7
$ws.Range("B10").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C10").Value = "Ñuble"
$ws.Range("D10").Value = 45264
$ws.Range("D10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E10").Value = 16
$ws.Range("F10").Value = "Fruta"
$ws.Range("G10").Value = 100103
$ws.Range("H10").Value = "Frutos de hueso (carozo)"
$ws.Range("I10").Value = 100103003
$ws.Range("J10").Value = "Damasco"
$ws.Range("K10").Value = "Dina"
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 24000
$ws.Range("O10").Value = 24000
$ws.Range("P10").Value = 24000
$ws.Range("Q10").Value = "`$/caja 10 kilos"
$ws.Range("R10").Value = "Provincia de Limarí"
$ws.Range("S10").Value = 2400
$ws.Range("T10").Value = 10

# Insert another new data row before the current row 17 (old row16, now at 17, -> new row18)
$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value = 7
$ws.Range("B17").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C17").Value = "Ñuble"
$ws.Range("D17").Value = 45265
$ws.Range("D17").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E17").Value = 16
$ws.Range("F17").Value = "Fruta"
$ws.Range("G17").Value = 100103
$ws.Range("H17").Value = "Frutos de hueso (carozo)"
$ws.Range("I17").Value = 100103003
$ws.Range("J17").Value = "Damasco"
$ws.Range("K17").Value = "Dina"
$ws.Range("L17").Value = "Primera"
$ws.Range("M17").Value = 80
$ws.Range("N17").Value = 20000
$ws.Range("O17").Value = 20000
$ws.Range("P17").Value = 20000
$ws.Range("Q17").Value = "`$/caja 10 kilos"
$ws.Range("R17").Value = "Provincia de Limarí"
$ws.Range("S17").Value = 2000
$ws.Range("T17").Value = 10

Write-Host "Rows now: $($ws.UsedRange.Rows.Count)"
